$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking text values (prices) must stay text: format cells as Text
# before assignment, then restore the Normal style so no stray formatting
# is introduced (matches original un-styled cells).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '30.680.52'
Set-TextValue $ws.Range("E2") '  +0.82%  '
Set-TextValue $ws.Range("D3") '2.102.75'
Set-TextValue $ws.Range("D4") '1.002'
Set-TextValue $ws.Range("E4") '  +0.23%  '
Set-TextValue $ws.Range("D5") '330.03'
Set-TextValue $ws.Range("E5") '  +1.91%  '
Set-TextValue $ws.Range("D6") '1.001'
Set-TextValue $ws.Range("E6") '  +0.20%  '
Set-TextValue $ws.Range("D7") '0.5270'
Set-TextValue $ws.Range("E7") '  +3.20%  '
Set-TextValue $ws.Range("D8") '0.4330'
Set-TextValue $ws.Range("E8") '  +4.79%  '
Set-TextValue $ws.Range("D9") '0.08924'
Set-TextValue $ws.Range("E9") '  +2.41%  '
Set-TextValue $ws.Range("D10") '46.64'
Set-TextValue $ws.Range("E10") '  +9.20%  '
Set-TextValue $ws.Range("D11") '1.167'
Set-TextValue $ws.Range("E11") '  +2.70%  '
Set-TextValue $ws.Range("D12") '24.59'
Set-TextValue $ws.Range("E12") '  -0.72%  '
Set-TextValue $ws.Range("D13") '2.106.05'
Set-TextValue $ws.Range("E13") '  +5.52%  '
Set-TextValue $ws.Range("D14") '6.691'
Set-TextValue $ws.Range("E14") '  +2.25%  '
Set-TextValue $ws.Range("D15") '7.790'
Set-TextValue $ws.Range("E15") '  +4.76%  '
Set-TextValue $ws.Range("D16") '96.85'
Set-TextValue $ws.Range("E16") '  +2.93%  '
Set-TextValue $ws.Range("D17") '1.002'
Set-TextValue $ws.Range("E17") '  +0.26%  '
Set-TextValue $ws.Range("D18") '0.00001125'
Set-TextValue $ws.Range("E18") '  +0.93%  '
Set-TextValue $ws.Range("D19") '0.06656'
Set-TextValue $ws.Range("E19") '  +2.28%  '
Set-TextValue $ws.Range("D20") '18.93'
Set-TextValue $ws.Range("E20") '  -0.03%  '
Set-TextValue $ws.Range("D21") '1.001'
Set-TextValue $ws.Range("E21") '  +0.18%  '
Set-TextValue $ws.Range("D22") '6.293'
Set-TextValue $ws.Range("E22") '  +1.97%  '
Set-TextValue $ws.Range("D23") '30.744.00'
Set-TextValue $ws.Range("E23") '  +0.91%  '
Set-TextValue $ws.Range("D24") '12.31'
Set-TextValue $ws.Range("E24") '  +4.13%  '
Set-TextValue $ws.Range("D25") '2.348.65'
Set-TextValue $ws.Range("E25") '  +5.23%  '
Set-TextValue $ws.Range("D26") '2.286'
Set-TextValue $ws.Range("E26") '  +3.41%  '
Set-TextValue $ws.Range("D27") '22.47'
Set-TextValue $ws.Range("E27") '  -0.09%  '
Set-TextValue $ws.Range("D28") '2.563'
Set-TextValue $ws.Range("E28") '  +5.98%  '
Set-TextValue $ws.Range("D29") '161.84'
Set-TextValue $ws.Range("E29") '  -0.85%  '
Set-TextValue $ws.Range("D30") '132.83'
Set-TextValue $ws.Range("E30") '  +0.73%  '
Set-TextValue $ws.Range("D31") '1.205'
Set-TextValue $ws.Range("E31") '  +5.70%  '
Set-TextValue $ws.Range("D32") '0.1076'
Set-TextValue $ws.Range("E32") '  +2.29%  '
Set-TextValue $ws.Range("D33") '6.143'
Set-TextValue $ws.Range("E33") '  +1.25%  '
Set-TextValue $ws.Range("D34") '1.548'
Set-TextValue $ws.Range("E34") '  +16.11%  '
Set-TextValue $ws.Range("D35") '3.854'
Set-TextValue $ws.Range("E35") '  +0.49%  '
Set-TextValue $ws.Range("D36") '0.02588'
Set-TextValue $ws.Range("E36") '  +3.12%  '
Set-TextValue $ws.Range("D37") '9.671'
Set-TextValue $ws.Range("E37") '  +7.30%  '
Set-TextValue $ws.Range("D38") '5.526'
Set-TextValue $ws.Range("E38") '  +2.81%  '
Set-TextValue $ws.Range("D39") '0.06717'
Set-TextValue $ws.Range("E39") '  +1.78%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D40") '12.60'
Set-TextValue $ws.Range("E40") '  +3.32%  '
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D41") '0.2270'
Set-TextValue $ws.Range("E41") '  +3.26%  '
Set-TextValue $ws.Range("D42") '0.6792'
Set-TextValue $ws.Range("E42") '  +2.53%  '
Set-TextValue $ws.Range("E43") '  +1.66%  '
Set-TextValue $ws.Range("D44") '1.001'
Set-TextValue $ws.Range("E44") '  +0.16%  '
Set-TextValue $ws.Range("D45") '0.6384'
Set-TextValue $ws.Range("E45") '  +3.74%  '
Set-TextValue $ws.Range("D46") '13.98'
Set-TextValue $ws.Range("E46") '  +2.33%  '
Set-TextValue $ws.Range("D47") '2.212'
Set-TextValue $ws.Range("E47") '  +0.20%  '
Set-TextValue $ws.Range("E48") '  -0.94%  '
Set-TextValue $ws.Range("D49") '1.252'
Set-TextValue $ws.Range("E49") '  -0.47%  '
Set-TextValue $ws.Range("D50") '82.74'
Set-TextValue $ws.Range("E50") '  +2.94%  '
Set-TextValue $ws.Range("E51") '  +6.56%  '
